# The workbook originally has a single sheet ("Sheet1") holding a small
# invoice-status table (Invoice Number / Items / Timestamp / Status).
# This "Update - Most Recent" edit duplicates that sheet (same layout,
# same header/status text, same styles) to capture a newer snapshot,
# naming the new tab after the source PDF it was generated from, and
# bumps the row-3 Timestamp to the newer value.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate Sheet1 (keeps formatting/styles/shared strings) and place
# the copy immediately after it.
$ws1.Copy($null, $ws1)
$newWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs.Name = "Sheet1 Sample Invoice_32.pdf"

# Refresh the "most recent" timestamp on the duplicated sheet.
$newWs.Cells.Item(3, 3).Value = 44945.5892592593

# Keep the original sheet as the active/selected tab.
$ws1.Activate()
